# Apply updated cryptocurrency price/volume snapshot values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is Cell -> NewValue. Updates are applied to the "cryptos" sheet
# (columns B/C for a couple of re-ranked coins, D = Price, E = Volume(1h)).
$updates = @(
    @{ Cell = "D2"; Value = "58.869.47" }
    @{ Cell = "E2"; Value = "  -6.41%  " }
    @{ Cell = "D3"; Value = "2.440.33" }
    @{ Cell = "E3"; Value = "  -9.13%  " }
    @{ Cell = "E4"; Value = "  +0.01%  " }
    @{ Cell = "D5"; Value = "536.99" }
    @{ Cell = "E5"; Value = "  -3.19%  " }
    @{ Cell = "D6"; Value = "146.29" }
    @{ Cell = "E6"; Value = "  -7.45%  " }
    @{ Cell = "E7"; Value = "  -0.14%  " }
    @{ Cell = "D8"; Value = "0.573" }
    @{ Cell = "E8"; Value = "  -2.67%  " }
    @{ Cell = "D9"; Value = "2.455.06" }
    @{ Cell = "E9"; Value = "  -8.70%  " }
    @{ Cell = "D10"; Value = "0.0989" }
    @{ Cell = "E10"; Value = "  -6.74%  " }
    @{ Cell = "E11"; Value = "  -2.23%  " }
    @{ Cell = "E12"; Value = "  -1.40%  " }
    @{ Cell = "E13"; Value = "  -5.09%  " }
    @{ Cell = "D14"; Value = "2.879.37" }
    @{ Cell = "E14"; Value = "  -8.85%  " }
    @{ Cell = "D15"; Value = "23.87" }
    @{ Cell = "E15"; Value = "  -10.27%  " }
    @{ Cell = "D16"; Value = "58.775.31" }
    @{ Cell = "E16"; Value = "  -6.44%  " }
    @{ Cell = "E17"; Value = "  -6.47%  " }
    @{ Cell = "D18"; Value = "2.501.69" }
    @{ Cell = "E18"; Value = "  -6.80%  " }
    @{ Cell = "D19"; Value = "11.11" }
    @{ Cell = "E19"; Value = "  -6.47%  " }
    @{ Cell = "E20"; Value = "  -5.90%  " }
    @{ Cell = "D21"; Value = "323.38" }
    @{ Cell = "E21"; Value = "  -6.37%  " }
    @{ Cell = "E22"; Value = "  -3.39%  " }
    @{ Cell = "D23"; Value = "5.68" }
    @{ Cell = "E23"; Value = "  -8.82%  " }
    @{ Cell = "D24"; Value = "60.58" }
    @{ Cell = "E24"; Value = "  -4.16%  " }
    @{ Cell = "E25"; Value = "  -11.32%  " }
    @{ Cell = "E26"; Value = "  -5.47%  " }
    @{ Cell = "E27"; Value = "  -2.32%  " }
    @{ Cell = "D28"; Value = "7.64" }
    @{ Cell = "E28"; Value = "  -6.62%  " }
    @{ Cell = "D29"; Value = "1.81" }
    @{ Cell = "E29"; Value = "  -6.54%  " }
    @{ Cell = "E30"; Value = "  -10.59%  " }
    @{ Cell = "D31"; Value = "6.63" }
    @{ Cell = "E31"; Value = "  -8.72%  " }
    @{ Cell = "D32"; Value = "1.23" }
    @{ Cell = "E32"; Value = "  -14.19%  " }
    @{ Cell = "E33"; Value = "  -0.05%  " }
    @{ Cell = "D34"; Value = "156.50" }
    @{ Cell = "E34"; Value = "  -4.66%  " }
    @{ Cell = "B35"; Value = "EthereumClassic" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" }
    @{ Cell = "D35"; Value = "18.51" }
    @{ Cell = "E35"; Value = "  -5.08%  " }
    @{ Cell = "B36"; Value = "ImmutableX" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D36"; Value = "1.37" }
    @{ Cell = "E36"; Value = "  -7.35%  " }
    @{ Cell = "D37"; Value = "4.43" }
    @{ Cell = "E37"; Value = "  -9.97%  " }
    @{ Cell = "D38"; Value = "1.68" }
    @{ Cell = "E38"; Value = "  -6.10%  " }
    @{ Cell = "D39"; Value = "5.82" }
    @{ Cell = "E39"; Value = "  -7.08%  " }
    @{ Cell = "D40"; Value = "311.90" }
    @{ Cell = "E40"; Value = "  -10.55%  " }
    @{ Cell = "E41"; Value = "  -5.72%  " }
    @{ Cell = "D42"; Value = "0.831" }
    @{ Cell = "E42"; Value = "  -12.25%  " }
    @{ Cell = "D43"; Value = "3.69" }
    @{ Cell = "E43"; Value = "  -7.54%  " }
    @{ Cell = "D44"; Value = "0.998" }
    @{ Cell = "E44"; Value = "  -0.08%  " }
    @{ Cell = "D45"; Value = "10.72" }
    @{ Cell = "E45"; Value = "  -2.43%  " }
    @{ Cell = "D46"; Value = "0.0940" }
    @{ Cell = "E46"; Value = "  -3.12%  " }
    @{ Cell = "E47"; Value = "  -5.96%  " }
    @{ Cell = "E48"; Value = "  -6.02%  " }
    @{ Cell = "E49"; Value = "  -5.28%  " }
    @{ Cell = "D50"; Value = "121.63" }
    @{ Cell = "E50"; Value = "  -5.52%  " }
    @{ Cell = "B51"; Value = "EnergySwap" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D51"; Value = "18.38" }
    @{ Cell = "E51"; Value = "  -8.89%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Price-looking numeric strings (e.g. "536.99") must stay text so Excel does not
    # silently coerce them into floating point numbers and drop trailing zeros.
    if ($u.Cell -match "^D" -and $u.Value -match "^\d+\.\d+$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
